# Auto-generated edit script applying F-column ("想去人数") updates
# and the single G8 type/value change on sheet "展览" (row 8: "暂时售罄" -> 0)
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 301
$ws.Range("F3").Value = 995
$ws.Range("F4").Value = 102
$ws.Range("F5").Value = 0
$ws.Range("F8").Value = 39592
$ws.Range("F9").Value = 0
$ws.Range("F11").Value = 8416
$ws.Range("F12").Value = 158
$ws.Range("F13").Value = 540
$ws.Range("F14").Value = 698
$ws.Range("F15").Value = 562
$ws.Range("F16").Value = 110
$ws.Range("F17").Value = 212
$ws.Range("F18").Value = 674
$ws.Range("F19").Value = 52
$ws.Range("F20").Value = 90
$ws.Range("F21").Value = 565
$ws.Range("F22").Value = 226
$ws.Range("F24").Value = 348
$ws.Range("F25").Value = 557
$ws.Range("F30").Value = 23
$ws.Range("F31").Value = 16
$ws.Range("F32").Value = 0
$ws.Range("F34").Value = 3
$ws.Range("F35").Value = 139
$ws.Range("F36").Value = 832
$ws.Range("F37").Value = 370
$ws.Range("F38").Value = 12
$ws.Range("F39").Value = 175
$ws.Range("F40").Value = 0
$ws.Range("F41").Value = 247
$ws.Range("F42").Value = 1010
$ws.Range("F44").Value = 1039
$ws.Range("F45").Value = 335
$ws.Range("F47").Value = 0
$ws.Range("F48").Value = 15

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 325
$ws.Range("F6").Value = 4382
$ws.Range("F8").Value = 301
$ws.Range("F9").Value = 7
$ws.Range("F10").Value = 5
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 6
$ws.Range("F14").Value = 0
$ws.Range("F16").Value = 7
$ws.Range("F17").Value = 72
$ws.Range("F18").Value = 0
$ws.Range("F20").Value = 4361
$ws.Range("F21").Value = 14
$ws.Range("F22").Value = 5

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 1783
$ws.Range("F3").Value = 425
$ws.Range("F4").Value = 329
$ws.Range("F5").Value = 0

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1783
$ws.Range("F4").Value = 329
$ws.Range("F5").Value = 301
$ws.Range("F6").Value = 102
$ws.Range("F7").Value = 943
$ws.Range("F8").Value = 1533
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 193
$ws.Range("F12").Value = 7
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 8416
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 698
$ws.Range("F20").Value = 562
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 92
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 674
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 90
$ws.Range("F28").Value = 565
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 1067
$ws.Range("F31").Value = 557
$ws.Range("F32").Value = 388
$ws.Range("F33").Value = 560
$ws.Range("F34").Value = 583
$ws.Range("F36").Value = 23
$ws.Range("F37").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("F40").Value = 832
$ws.Range("F41").Value = 370
$ws.Range("F42").Value = 175
$ws.Range("F43").Value = 247
$ws.Range("F44").Value = 211
$ws.Range("F45").Value = 1039
$ws.Range("F46").Value = 335
$ws.Range("F47").Value = 0
$ws.Range("F48").Value = 15
$ws.Range("F49").Value = 5

# Special case: G8 on sheet "展览" changes from text "暂时售罄" to numeric 0
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("G8").Value = 0
